# core: update agent prompt
#
# Add a "Color Indicator" header label in the first empty column (G1)
# to the right of the existing x/y1/y2 table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Color Indicator"
